# Updated cryptos list with refreshed Price (column D) and Volume(1h) (column E)
# values, mirroring the upstream GitHub Actions data refresh.
#
# Column D values are written with a leading apostrophe so Excel keeps them
# as text (matching the original sheet, where these "price" strings such as
# "27.716.24" or "1.001" are stored as plain text, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.727.12'
$ws.Range("D3").Value = '''1.901.32'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''311.83'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '''0.5231'
$ws.Range("E7").Value = '  +5.97%  '
$ws.Range("D8").Value = '''0.3790'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.07248'
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("D10").Value = '''21.24'
$ws.Range("E10").Value = '  +3.07%  '
$ws.Range("D11").Value = '''0.9029'
$ws.Range("D12").Value = '''0.07644'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '''1.896.99'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '''5.446'
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").Value = '''92.31'
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '''0.000008694'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '''27.760.72'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = '''14.47'
$ws.Range("D21").Value = '''5.140'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").Value = '''2.157.85'
$ws.Range("E22").Value = '  +0.91%  '
$ws.Range("E23").Value = '  +0.90%  '
$ws.Range("D24").Value = '''6.618'
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").Value = '''153.12'
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("D26").Value = '''1.868'
$ws.Range("E26").Value = '  +0.96%  '
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").Value = '''2.162'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D29").Value = '''114.48'
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("D30").Value = '''4.839'
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").Value = '''0.09078'
$ws.Range("E31").Value = '  +1.56%  '
$ws.Range("D32").Value = '''3.188'
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("D33").Value = '''4.836'
$ws.Range("E33").Value = '  +4.05%  '
$ws.Range("D34").Value = '''1.224'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("D35").Value = '''0.7788'
$ws.Range("D36").Value = '''0.02092'
$ws.Range("E36").Value = '  +2.25%  '
$ws.Range("D37").Value = '''2.575'
$ws.Range("E38").Value = '  +2.69%  '
$ws.Range("D39").Value = '''1.093'
$ws.Range("E39").Value = '  -0.63%  '
$ws.Range("D40").Value = '''0.5547'
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("D41").Value = '''0.05287'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '''6.733'
$ws.Range("E42").Value = '  -2.48%  '
$ws.Range("D43").Value = '''116.24'
$ws.Range("D44").Value = '''8.512'
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").Value = '''0.1518'
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").Value = '''0.4816'
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").Value = '''10.51'
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("D48").Value = '''0.9997'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").Value = '''66.79'
$ws.Range("E50").Value = '  -0.98%  '
$ws.Range("D51").Value = '''0.06002'
$ws.Range("E51").Value = '  -0.96%  '
